# optimise the code of linkgame
#
# The "mini" sheet lists minigame buttons. The row for 连连看 (LinkGame),
# row 10, had its IconPath (column D) pointing at the "GameButton6" icon
# (the same icon reused by row 9 / 俄罗斯块). Give LinkGame its own icon
# name, "GameButton7", and move the active selection to D8.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mini")

# Row 10 = Id 17000007 / 连连看 -> change IconPath from "GameButton6" to
# the new unique shared string "GameButton7".
$ws.Range("D10").Value = "GameButton7"

# Move the selected/active cell to D8, as saved in the sheet view.
$ws.Range("D8").Select()
